$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old row 14 ("DIO:" block), pushing it down to rows 16-19
$ws.Rows("14:15").Insert()

# New motor map entries (order matches shared-string insertion order)
$ws.Range("B11").Value = "climber master"
$ws.Range("B12").Value = "lifter"
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "climber slave"
$ws.Range("B6").Value = "shooter"
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = "collector stirrer"

# Leave the active cell/selection on B14, matching the saved view state
$ws.Range("B14").Select()
